$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-40: Fecha, Calidad, Volumen, Precio minimo/maximo/promedio, Precio $/Kg ---

$ws.Range("D2").Value = 44487
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 23000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 23500
$ws.Range("S2").Value = 2350

$ws.Range("D3").Value = 45194
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 22000
$ws.Range("S3").Value = 2200

$ws.Range("D4").Value = 44868
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 26000
$ws.Range("O4").Value = 26000
$ws.Range("P4").Value = 26000
$ws.Range("S4").Value = 2600

$ws.Range("D5").Value = 45247
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 2000

$ws.Range("D6").Value = 45236
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 22000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 22000
$ws.Range("S6").Value = 2200

$ws.Range("D7").Value = 45216
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 21000
$ws.Range("P7").Value = 21000
$ws.Range("S7").Value = 2100

$ws.Range("D8").Value = 45212
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 22000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 22000
$ws.Range("S8").Value = 2200

$ws.Range("D9").Value = 44446
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 21000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 21500
$ws.Range("S9").Value = 2150

$ws.Range("D10").Value = 45209
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 22000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 22000
$ws.Range("S10").Value = 2200

$ws.Range("D11").Value = 44841
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 23000
$ws.Range("O11").Value = 24000
$ws.Range("P11").Value = 23500
$ws.Range("S11").Value = 2350

$ws.Range("D12").Value = 45176
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 22000
$ws.Range("O12").Value = 22000
$ws.Range("P12").Value = 22000
$ws.Range("S12").Value = 2200

$ws.Range("D13").Value = 45205
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 22000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 22000
$ws.Range("S13").Value = 2200

$ws.Range("D14").Value = 44460
$ws.Range("L14").Value = 'Especial'
$ws.Range("M14").Value = 60
$ws.Range("N14").Value = 31000
$ws.Range("O14").Value = 32000
$ws.Range("P14").Value = 31500
$ws.Range("S14").Value = 3150

$ws.Range("D15").Value = 44460
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 30000
$ws.Range("O15").Value = 30000
$ws.Range("P15").Value = 30000
$ws.Range("S15").Value = 3000

$ws.Range("D16").Value = 45191
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 21000
$ws.Range("O16").Value = 21000
$ws.Range("P16").Value = 21000
$ws.Range("S16").Value = 2100

$ws.Range("D17").Value = 44447
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 21000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 21500
$ws.Range("S17").Value = 2150

$ws.Range("D18").Value = 44461
$ws.Range("L18").Value = 'Especial'
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = 31000
$ws.Range("O18").Value = 32000
$ws.Range("P18").Value = 31500
$ws.Range("S18").Value = 3150

$ws.Range("D19").Value = 44461
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = 30000
$ws.Range("O19").Value = 30000
$ws.Range("P19").Value = 30000
$ws.Range("S19").Value = 3000

$ws.Range("D20").Value = 44448
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = 21000
$ws.Range("O20").Value = 22000
$ws.Range("P20").Value = 21500
$ws.Range("S20").Value = 2150

$ws.Range("D21").Value = 45173
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 22000
$ws.Range("O21").Value = 22000
$ws.Range("P21").Value = 22000
$ws.Range("S21").Value = 2200

$ws.Range("D22").Value = 45203
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 30
$ws.Range("N22").Value = 21000
$ws.Range("O22").Value = 21000
$ws.Range("P22").Value = 21000
$ws.Range("S22").Value = 2100

$ws.Range("D23").Value = 44874
$ws.Range("L23").Value = 'Especial'
$ws.Range("M23").Value = 30
$ws.Range("N23").Value = 25000
$ws.Range("O23").Value = 25000
$ws.Range("P23").Value = 25000
$ws.Range("S23").Value = 2500

$ws.Range("D24").Value = 44874
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 80
$ws.Range("N24").Value = 23000
$ws.Range("O24").Value = 24000
$ws.Range("P24").Value = 23500
$ws.Range("S24").Value = 2350

$ws.Range("D25").Value = 45196
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 23000
$ws.Range("O25").Value = 23000
$ws.Range("P25").Value = 23000
$ws.Range("S25").Value = 2300

$ws.Range("D26").Value = 44452
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 60
$ws.Range("N26").Value = 21000
$ws.Range("O26").Value = 22000
$ws.Range("P26").Value = 21500
$ws.Range("S26").Value = 2150

$ws.Range("D27").Value = 45230
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 21000
$ws.Range("O27").Value = 21000
$ws.Range("P27").Value = 21000
$ws.Range("S27").Value = 2100

$ws.Range("D28").Value = 45189
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 22000
$ws.Range("O28").Value = 22000
$ws.Range("P28").Value = 22000
$ws.Range("S28").Value = 2200

$ws.Range("D29").Value = 45219
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 30
$ws.Range("N29").Value = 20000
$ws.Range("O29").Value = 20000
$ws.Range("P29").Value = 20000
$ws.Range("S29").Value = 2000

$ws.Range("D30").Value = 45224
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 80
$ws.Range("N30").Value = 20000
$ws.Range("O30").Value = 20000
$ws.Range("P30").Value = 20000
$ws.Range("S30").Value = 2000

$ws.Range("D31").Value = 44839
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 120
$ws.Range("N31").Value = 25000
$ws.Range("O31").Value = 26000
$ws.Range("P31").Value = 25500
$ws.Range("S31").Value = 2550

$ws.Range("D32").Value = 45239
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 100
$ws.Range("N32").Value = 22000
$ws.Range("O32").Value = 22000
$ws.Range("P32").Value = 22000
$ws.Range("S32").Value = 2200

$ws.Range("D33").Value = 45237
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 80
$ws.Range("N33").Value = 22000
$ws.Range("O33").Value = 22000
$ws.Range("P33").Value = 22000
$ws.Range("S33").Value = 2200

$ws.Range("D34").Value = 45237
$ws.Range("L34").Value = 'Segunda'
$ws.Range("M34").Value = 50
$ws.Range("N34").Value = 18000
$ws.Range("O34").Value = 18000
$ws.Range("P34").Value = 18000
$ws.Range("S34").Value = 1800

$ws.Range("D35").Value = 45225
$ws.Range("L35").Value = 'Primera'
$ws.Range("M35").Value = 80
$ws.Range("N35").Value = 21000
$ws.Range("O35").Value = 21000
$ws.Range("P35").Value = 21000
$ws.Range("S35").Value = 2100

$ws.Range("D36").Value = 45217
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 30
$ws.Range("N36").Value = 21000
$ws.Range("O36").Value = 21000
$ws.Range("P36").Value = 21000
$ws.Range("S36").Value = 2100

$ws.Range("D37").Value = 45244
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 90
$ws.Range("N37").Value = 20000
$ws.Range("O37").Value = 21000
$ws.Range("P37").Value = 20444
$ws.Range("S37").Value = 2044

$ws.Range("D38").Value = 44848
$ws.Range("L38").Value = 'Especial'
$ws.Range("M38").Value = 60
$ws.Range("N38").Value = 24000
$ws.Range("O38").Value = 25000
$ws.Range("P38").Value = 24500
$ws.Range("S38").Value = 2450

$ws.Range("D39").Value = 44848
$ws.Range("L39").Value = 'Primera'
$ws.Range("M39").Value = 120
$ws.Range("N39").Value = 21000
$ws.Range("O39").Value = 22000
$ws.Range("P39").Value = 21500
$ws.Range("S39").Value = 2150

$ws.Range("D40").Value = 45243
$ws.Range("L40").Value = 'Primera'
$ws.Range("M40").Value = 60
$ws.Range("N40").Value = 21000
$ws.Range("O40").Value = 21000
$ws.Range("P40").Value = 21000
$ws.Range("S40").Value = 2100

# --- Append new data row 41 ---
$ws.Range("A41").Value = 7
$ws.Range("B41").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C41").Value = 'Ñuble'
$ws.Range("D41").Value = 45180
$ws.Range("D41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E41").Value = 16
$ws.Range("F41").Value = 'Fruta'
$ws.Range("G41").Value = 100107
$ws.Range("H41").Value = 'Otros'
$ws.Range("I41").Value = 100107002
$ws.Range("J41").Value = 'Chirimoya'
$ws.Range("K41").Value = 'Cultivar IV Región'
$ws.Range("L41").Value = 'Primera'
$ws.Range("M41").Value = 40
$ws.Range("N41").Value = 22000
$ws.Range("O41").Value = 22000
$ws.Range("P41").Value = 22000
$ws.Range("Q41").Value = '$/bandeja 10 kilos'
$ws.Range("R41").Value = 'Provincia de Limarí'
$ws.Range("S41").Value = 2200
$ws.Range("T41").Value = 10
